$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 205 (this shifts old rows 205-208 down to 206-209,
# carrying their existing values/formatting along).
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record.
$ws.Cells.Item(205, 1).Value = 5
$ws.Cells.Item(205, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(205, 3).Value = "Maule"
$ws.Cells.Item(205, 4).Value = 44628
$ws.Cells.Item(205, 5).Value = 7
$ws.Cells.Item(205, 6).Value = 100112024
$ws.Cells.Item(205, 7).Value = "Choclo"
$ws.Cells.Item(205, 8).Value = "Choclero"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 50000
$ws.Cells.Item(205, 11).Value = 130
$ws.Cells.Item(205, 12).Value = 150
$ws.Cells.Item(205, 13).Value = 142
$ws.Cells.Item(205, 14).Value = "$/unidad"
$ws.Cells.Item(205, 15).Value = "Región del Maule"
$ws.Cells.Item(205, 16).Value = 142
$ws.Cells.Item(205, 17).Value = 1
$ws.Cells.Item(205, 18).Value = "Hortaliza"
